$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 7; $r -le 14; $r++) {
    $ws.Cells.Item($r, 1).Value = "A$r"
    $ws.Cells.Item($r, 2).Value = "B$r"
    $ws.Cells.Item($r, 3).Value = "C$r"
}
